$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item("Tableau1")
$newRow = $tbl.ListRows.Add()

$ws.Range("E53:M53").Copy()
$ws.Range("E54:M54").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Rows.Item(54).RowHeight = 43.2

$ws.Range("E54").Value = 44287
$ws.Range("F54").Value = 0.79166666666666663
$ws.Range("G54").Value = 0.80208333333333337
$ws.Range("H54").Formula = "=IF(ISBLANK(Tableau1[[#This Row],[Heure Début]]),"""",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure Début]])"
$ws.Range("I54").Value = "Développement"
$ws.Range("J54").Value = "Corriger le décalage avec la grille"
$ws.Range("K54").Value = "Domicile"
$ws.Range("L54").Value = "Correction du décalage de 1 avec la grille"

$ws.Range("L55").Select() | Out-Null
